$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '28.547.55'
Set-TextValue "E2" '  -0.11%  '

Set-TextValue "D3" '1.825.52'
Set-TextValue "E3" '  -0.13%  '

Set-TextValue "E4" '  +0.04%  '

Set-TextValue "D5" '315.69'

Set-TextValue "E6" '  +0.03%  '

Set-TextValue "D7" '0.5115'
Set-TextValue "E7" '  -5.57%  '

Set-TextValue "E8" '  -1.56%  '

Set-TextValue "D9" '0.08220'
Set-TextValue "E9" '  +7.12%  '

Set-TextValue "E10" '  -0.66%  '

Set-TextValue "D11" '41.70'
Set-TextValue "E11" '  -0.35%  '

Set-TextValue "D12" '21.21'
Set-TextValue "E12" '  +0.25%  '

Set-TextValue "D13" '6.337'

Set-TextValue "E14" '  +0.03%  '

Set-TextValue "D15" '7.539'
Set-TextValue "E15" '  -1.35%  '

Set-TextValue "D16" '1.834.38'
Set-TextValue "E16" '  +0.43%  '

Set-TextValue "D17" '0.00001130'
Set-TextValue "E17" '  +3.58%  '

Set-TextValue "D18" '92.85'
Set-TextValue "E18" '  +3.12%  '

Set-TextValue "E19" '  +0.97%  '

Set-TextValue "D20" '17.85'

Set-TextValue "D21" '1.000'
Set-TextValue "E21" '  +0.02%  '

Set-TextValue "D22" '6.105'
Set-TextValue "E22" '  +0.59%  '

Set-TextValue "D23" '28.580.53'
Set-TextValue "E23" '  -0.03%  '

Set-TextValue "D24" '11.42'
Set-TextValue "E24" '  +2.06%  '

Set-TextValue "D25" '2.265'
Set-TextValue "E25" '  -0.48%  '

Set-TextValue "D26" '21.43'
Set-TextValue "E26" '  +3.24%  '

Set-TextValue "E27" '  -0.82%  '

Set-TextValue "D28" '2.036.09'
Set-TextValue "E28" '  +0.01%  '

Set-TextValue "D29" '2.410'
Set-TextValue "E29" '  -1.97%  '

Set-TextValue "D30" '126.96'
Set-TextValue "E30" '  +2.30%  '

Set-TextValue "D31" '1.114'
Set-TextValue "E31" '  -1.50%  '

Set-TextValue "D32" '0.1091'
Set-TextValue "E32" '  -1.80%  '

Set-TextValue "D33" '5.780'
Set-TextValue "E33" '  +1.62%  '

Set-TextValue "D34" '3.657'
Set-TextValue "E34" '  +0.29%  '

Set-TextValue "E35" '  -6.10%  '

Set-TextValue "D36" '0.2230'
Set-TextValue "E36" '  -0.84%  '

Set-TextValue "D37" '5.286'
Set-TextValue "E37" '  +1.37%  '

Set-TextValue "D38" '0.02358'
Set-TextValue "E38" '  -0.22%  '

Set-TextValue "D39" '8.853'
Set-TextValue "E39" '  -0.37%  '

Set-TextValue "D40" '0.6329'
Set-TextValue "E40" '  +0.39%  '

Set-TextValue "E41" '  -0.54%  '

Set-TextValue "D42" '1.183'
Set-TextValue "E42" '  -0.48%  '

Set-TextValue "B43" 'WEMIXTOKEN'
Set-TextValue "C43" 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D43" '1.398'
Set-TextValue "E43" '  -0.31%  '

Set-TextValue "B44" 'EnergySwap'
Set-TextValue "C44" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D44" '13.61'
Set-TextValue "E44" '  +1.16%  '

Set-TextValue "B45" 'Decentraland'
Set-TextValue "C45" 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue "D45" '0.5957'
Set-TextValue "E45" '  +1.14%  '

Set-TextValue "B46" 'PancakeSwap'
Set-TextValue "C46" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D46" '3.736'
Set-TextValue "E46" '  +0.67%  '

Set-TextValue "B47" 'Quant'
Set-TextValue "C47" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D47" '125.72'
Set-TextValue "E47" '  +0.46%  '

Set-TextValue "B48" 'NEARProtocol'
Set-TextValue "C48" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D48" '1.997'
Set-TextValue "E48" '  -0.46%  '

Set-TextValue "B49" 'EOS'
Set-TextValue "C49" 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue "D49" '1.194'
Set-TextValue "E49" '  -0.44%  '

Set-TextValue "B50" 'Cronos'
Set-TextValue "C50" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D50" '0.06930'
Set-TextValue "E50" '  +0.28%  '

Set-TextValue "B51" 'ThetaToken'
Set-TextValue "C51" 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue "D51" '1.086'
Set-TextValue "E51" '  +4.73%  '
